# Weekly update: insert a new price record as row 3, pushing the
# existing history down by one row (matches the "Fruta / hortaliza,
# semanal" commit that adds the latest week's observation).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing rows 3:20 down to 4:21, inserting a blank row 3
# (format is inherited from the row above, same as Excel's native
# Insert behaviour).
$ws.Rows.Item(3).Insert()

# Populate the new row 3 with this week's observation.
$ws.Range("A3").Value2 = 1
$ws.Range("B3").Value2 = "Agrícola del Norte S.A. de Arica"
$ws.Range("C3").Value2 = "Arica y Parinacota"
$ws.Range("D3").Value2 = 44972
$ws.Range("E3").Value2 = 15
$ws.Range("F3").Value2 = 100114007
$ws.Range("G3").Value2 = "Jengibre"
$ws.Range("H3").Value2 = "Sin especificar"
$ws.Range("I3").Value2 = "Primera"
$ws.Range("J3").Value2 = 350
$ws.Range("K3").Value2 = 17000
$ws.Range("L3").Value2 = 18000
$ws.Range("M3").Value2 = 17429
$ws.Range("N3").Value2 = '$/caja 15 kilos'
$ws.Range("O3").Value2 = "Perú"
$ws.Range("P3").Value2 = 1162
$ws.Range("Q3").Value2 = 15
$ws.Range("R3").Value2 = "Hortaliza"
